# Horarios actualizados Linea 141 - 627
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173)
# with the latest scrape data (Ultima actualizacion 04:56:49).

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')

$ws.Range('A1').Value = 'LÍNEA 141 - LP1912 - 19/01/2026'
$ws.Range('A2').Value = 'Última actualización: 04:56:49'
$ws.Range('A3').Value = 'Total filas: 25'

$data = @(
    @('03:45:25', '03:47', '14_ABASTO', 2, 'LP1912'),
    @('03:45:25', '04:01', '81_EL PELIGRO', 16, 'LP1912'),
    @('03:45:25', '04:46', '215A_EL PATO', 61, 'LP1912'),
    @('03:45:25', '04:53', '11_ETCHEVERRY', 68, 'LP1912'),
    @('04:56:49', '05:13', '14_ABASTO', 17, 'LP1912'),
    @('03:45:25', '05:16', '17_ROMERO', 91, 'LP1912'),
    @('04:45:05', '05:16', '14_ABASTO', 31, 'LP1912'),
    @('03:45:25', '05:22', '23_HERNANDEZ', 97, 'LP1912'),
    @('03:45:25', '05:34', '215B_EL PATO', 109, 'LP1912'),
    @('04:18:02', '05:34', '14_ABASTO', 76, 'LP1912'),
    @('04:18:02', '05:35', '215B_EL PATO', 77, 'LP1912'),
    @('03:45:25', '05:37', '14_ABASTO', 112, 'LP1912'),
    @('04:18:02', '05:46', '15_ABASTO', 88, 'LP1912'),
    @('04:45:05', '06:04', '16_SANTA ANA', 79, 'LP1912'),
    @('04:18:02', '06:05', '16_SANTA ANA', 107, 'LP1912'),
    @('04:56:49', '06:11', '215A_EL PATO', 75, 'LP1912'),
    @('04:18:02', '06:12', '215A_EL PATO', 114, 'LP1912'),
    @('04:18:02', '06:14', '225_HARAS DEL SUR', 116, 'LP1912'),
    @('04:45:05', '06:21', '26_HERNANDEZ', 96, 'LP1912'),
    @('04:45:05', '06:27', '23_HERNANDEZ', 102, 'LP1912'),
    @('04:56:49', '06:29', '86_EST CHICA-ESC AGRARIA', 93, 'LP1912'),
    @('04:45:05', '06:30', '86_EST CHICA-ESC AGRARIA', 105, 'LP1912'),
    @('04:45:05', '06:31', '16_SANTA ANA', 106, 'LP1912'),
    @('04:45:05', '06:44', '225_C ROCA-H SUR', 119, 'LP1912'),
    @('04:56:49', '06:46', '215C_EL PATO', 110, 'LP1912')
)

$r = 6
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')

$ws.Range('A1').Value = 'LÍNEA 141 - LP1912-215 - 19/01/2026'
$ws.Range('A2').Value = 'Última actualización: 04:56:49'
$ws.Range('A3').Value = 'Total filas: 6'

$data = @(
    @('03:45:25', '04:46', '215A_EL PATO', 61, 'LP1912'),
    @('03:45:25', '05:34', '215B_EL PATO', 109, 'LP1912'),
    @('04:18:02', '05:35', '215B_EL PATO', 77, 'LP1912'),
    @('04:56:49', '06:11', '215A_EL PATO', 75, 'LP1912'),
    @('04:18:02', '06:12', '215A_EL PATO', 114, 'LP1912'),
    @('04:56:49', '06:46', '215C_EL PATO', 110, 'LP1912')
)

$r = 6
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')

$ws.Range('A1').Value = 'LÍNEA 141 - 6203-6173 - 19/01/2026'
$ws.Range('A2').Value = 'Última actualización: 04:56:49'
$ws.Range('A3').Value = 'Total filas: 6'

$data = @(
    @('04:56:49', '05:43', '215A_LA PLATA', 47, 'L6173'),
    @('03:45:25', '05:44', '215A_LA PLATA', 119, 'L6173'),
    @('04:56:49', '06:08', '215A_LA PLATA', 72, 'L6173'),
    @('04:18:02', '06:09', '215A_LA PLATA', 111, 'L6173'),
    @('04:56:49', '06:32', '215C_LA PLATA', 96, 'L6203'),
    @('04:45:05', '06:33', '215C_LA PLATA', 108, 'L6203')
)

$r = 6
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

